# --- edit.ps1 ---------------------------------------------------------
# Reproduces the commit "added new testcase InsuranceRegisterTest and
# minor enhancements" against DDTdata.xlsx:
#   * Sheet3 gets populated with sample rows (1..8 in column A, "Hello"
#     in B4) and its selection moves to D12.
#   * A brand-new "Emails" worksheet is appended at the end of the
#     workbook, gets "Hello" in E5, its selection set to J10, and is
#     left as the active/selected sheet (as in the target workbook.xml,
#     where activeTab now points at the 4th tab).
#   * ProjectDDs (Sheet2) gets two new plain-text labels in A1/A2
#     ("word 1" / "word 2 ") and its selection moves to B3; it is no
#     longer the active tab.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Sheet3: sample data -----------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
for ($i = 1; $i -le 8; $i++) {
    $ws3.Cells.Item($i, 1).Value = $i
}
$ws3.Range("B4").Value = "Hello"
[void]$ws3.Range("D12").Select()

# ---- New "Emails" worksheet, appended after the last sheet -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEmails = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsEmails.Name = "Emails"
$wsEmails.Range("E5").Value = "Hello"
[void]$wsEmails.Range("J10").Select()

# ---- ProjectDDs: new labels ---------------------------------------------
$ws2 = $wb.Worksheets.Item("ProjectDDs")
$ws2.Range("A1").ClearFormats()
$ws2.Range("A1").Value = "word 1"
$ws2.Range("A2").ClearFormats()
$ws2.Range("A2").Value = "word 2 "
[void]$ws2.Range("B3").Select()

# ---- Leave "Emails" as the active sheet/tab -----------------------------
[void]$wsEmails.Activate()
